$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the crypto table
# with the latest scraped figures. All cells in these columns are plain
# text (some prices contain multiple "." thousands separators, e.g.
# "65.807.62", and the percentages keep their original two-space padding
# on both sides, e.g. "  +3.31%  "). For price cells whose new text would
# otherwise be auto-recognized by Excel as a genuine number (losing a
# trailing zero such as "5.90" -> 5.9), the cell is forced to Text format
# just long enough to accept the literal string, then restored to the
# default "Normal" style so no stray formatting is left behind.
$ws.Range("D2").Value = "65.807.62"
$ws.Range("E2").Value = "  +3.31%  "
$ws.Range("D3").Value = "2.664.08"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  +9.38%  "
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +16.07%  "
$ws.Range("D15").Value = "3.151.41"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "65.548.69"
$ws.Range("E16").Value = "  +3.06%  "
$ws.Range("D17").Value = "2.673.29"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "361.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.28%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.31%  "
$ws.Range("E26").Value = "  +17.34%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("E36").Value = "  +4.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "167.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("E45").Value = "  +9.25%  "
$ws.Range("E46").Value = "  +5.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.664"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0266"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.64%  "
